$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "alias"
$ws.Range("D1").Value = "classID"
$ws.Range("D1").Select()
